# Fruta / hortaliza, semanal
# Insert a new data row at sheet row 135, pushing existing rows 135-185
# down to 136-186, and populate the new row with its values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 135 (shifts 135..185 -> 136..186)
$ws.Rows.Item(135).Insert()

# Populate the newly inserted row 135 with the new record.
$ws.Range("A135").Value = 4
$ws.Range("B135").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C135").Value = "Los Lagos"
$ws.Range("D135").Value = 44489
$ws.Range("E135").Value = 10
$ws.Range("F135").Value = 100112045
$ws.Range("G135").Value = "Zapallo"
$ws.Range("H135").Value = "Paine"
$ws.Range("I135").Value = "1a (guarda)"
$ws.Range("J135").Value = 150
$ws.Range("K135").Value = 400
$ws.Range("L135").Value = 400
$ws.Range("M135").Value = 400
$ws.Range("N135").Value = '$/kilo (volumen en unidades)'
$ws.Range("O135").Value = "Región Metropolitana"
$ws.Range("P135").Value = 400
$ws.Range("Q135").Value = 1
$ws.Range("R135").Value = "Hortaliza"
